$wb = $excel.ActiveWorkbook

# --- Foglio1: selection/view changes (tab no longer selected, range A1:B14 selected) ---
$ws1 = $wb.Worksheets.Item("Foglio1")
$ws1.Range("A1:B14").Select() | Out-Null

# --- new worksheet "liste_europee", placed after "europee" ---
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "liste_europee"

$data = @(
    @("COALIZIONE", "LISTA"),
    @("SINISTRA", "PARTITO DEMOCRATICO"),
    @("SINISTRA", "ALLEANZA VERDI E SINISTRA"),
    @("SINISTRA", "MOVIMENTO 5 STELLE"),
    @("SINISTRA", "STATI UNITI D'EUROPA"),
    @("SINISTRA", "AZIONE"),
    @("DESTRA", "FORZA ITALIA"),
    @("DESTRA", "LEGA"),
    @("DESTRA", "FRATELLI D'ITALIA"),
    @("PTD", "PACE TERRA DIGNITA'"),
    @("L", "LIBERTA'"),
    @("AP", "ALTERNATIVA POPOLARE"),
    @("SVP", "SVP")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 1
    $ws3.Cells.Item($row, 1).Value = $data[$i][0]
    $ws3.Cells.Item($row, 2).Value = $data[$i][1]
}

# header row bold, matching Foglio1's header style
$ws3.Range("A1:B1").Font.Bold = $true

# best-fit-ish width for column A
$ws3.Columns.Item(1).ColumnWidth = 10.86

# final selection / active sheet
$ws3.Range("B4").Select() | Out-Null
